# MainBoard_BOM.xlsx update
# Commit: Updated BOM file (removed R49, C117 and fixed C116 and J2).
#
# Row 14 (Item 11): C116, C117 -> C116 only; part corrected to a 250V 0805 cap.
# Row 17 (Item 14): J2 RJ45 connector corrected to the Pulse J0011D21BNL jack.
# Row 33 (Item 30): R49 removed from the DNP/0-ohm reference list, qty -1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 14: C116, C117 -> C116 (fixed part number / package / price) ---
$ws.Range("B14").Value = "C116"
$ws.Range("H14").Value = "445-2277-1-ND"
$ws.Range("F14").Value = "C2012X7R2E102K"
$ws.Range("I14").Value = "CAP CER 1000PF 250V 10% X7R 0805"
$ws.Range("D14").Value = "C0805"
$ws.Range("J14").Value = 1
$ws.Range("K14").Value = 0.17

# --- Row 17: J2 connector fixed to the Pulse Electronics J0011D21BNL jack ---
$ws.Range("H17").Value = "553-1485-ND"
$ws.Range("F17").Value = "J0011D21BNL"
$ws.Range("I17").Value = "CONN PULSEJACK 1PORT 10/100B-TX"
$ws.Range("C17").Value = "CON-RJ45-J0011D21BNL"
$ws.Range("D17").Value = "CON-RJ45-J0011D21BNL"
$ws.Range("K17").Value = 7.1

# --- Row 33: remove R49 from the reference designator list, qty 12 -> 11 ---
$ws.Range("B33").Value = "R4, R11, R12, R13, R15, R16, R17, R18, R74, R75"
$ws.Range("J33").Value = 11

# --- Restore view: clear the frozen/scrolled top-left cell and select J34 ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 85
$ws.Range("J34").Select() | Out-Null
